# Adds the 4 new records (rows 115-118) to Sheet1, matching the commit
# "Atualização automática do arquivo Excel".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rows = @(
    @{ Row = 115; A = "Globo"; B = "RJ TV 2"; C = "Obras";   D = "2025-04-10T19:27"; E = "Negativo"; F = "UBS demolida em Guarus. Moradores reclamam que precisam se deslocar para ter atendimento médico. Unidade Básica de Saúde Alair Ferreira, no Parque Vicente Dias, foi demolida. Moradores das proximidades sentem falta da unidade, que recebeu muitos elogios sobre os atendimentos médicos que eram oferecidos. Entrevista com moradoras. Em fevereiro, prefeitura anunciou um Centro Especializado de Reabilitação, um investimento de cerca de R$ 7 milhões. Produção questionou sobre atendimento e obras. Em nota, prefeitura só respondeu sobre atendimento: moradores podem procurar o Centro de Saúde de Guarus. `n*matéria*" },
    @{ Row = 116; A = "Globo"; B = "RJ TV 2"; C = "Obras";   D = "2025-04-10T19:30"; E = "Negativo"; F = "Por whatsApp, moradora da Rua Maricá, no Parque Guarus, reclama que a rua começou a receber asfalto pouco antes das eleições. Apenas, um lado recebeu asfalto. A cada dia, a rua fica pior. A rua é extensa, movimentada e sai na BR-101. *com nota* da prefeitura" },
    @{ Row = 117; A = "Globo"; B = "RJ TV 2"; C = "Saúde";   D = "2025-04-10T19:33"; E = "Negativo"; F = "Paralisação no Pronto-Socorro Pediátrico. No Plantadores de Cana, apenas casos graves estão sendo atendidos. Profissionais estariam 3 meses sem receber salários. Pronto Socorro Pediátrico mantém emergência vermelha. Atendimento pediátrico também é feito no PU de Guarus. Sobre o pagamento e o atendimento de baixa complexidade na unidade, não houve resposta. " },
    @{ Row = 118; A = "Globo"; B = "RJ TV 2"; C = "Governo"; D = "2025-04-10T19:38"; E = "Neutro";   F = "O que os prefeitos já fizeram pela sua cidade? Acompanhe nos nossos telejornais e no G1, a partir de segunda, um balanço dos 100 dias de governo.  " }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
}

# Multi-line text (e.g. the embedded line break in F115) makes the engine
# stamp an explicit custom row height on save; AutoFit restores the default
# (no ht / customHeight attributes) so the new rows look like the rest of
# the sheet.
$ws.Range("A115:A118").EntireRow.AutoFit() | Out-Null
